# Working user data input via number keyboard
#
# Updates the "Translation" sheet of the TouchGFX texts workbook:
#   - Rows 24 and 26 (SingleUseId30 / SingleUseId33) lose their unit
#     suffix, becoming plain "<value>" templates (the unit is now a
#     separate text so it can be composed with the numeric keyboard
#     input screen).
#   - Four new text rows are appended (32-35): SingleUseId39..42, adding
#     a standalone "<value> l" template, a new prompt
#     ("Twoje zapotrzebowanie na wodę wynosi:"), and standalone "cm"/"l"
#     unit strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# --- modify existing rows -------------------------------------------------

# SingleUseId30: "<value> cm" -> "<value>"
$ws.Range("F24").Value = "<value>"

# SingleUseId33: "<value> l" -> "<value>"
$ws.Range("F26").Value = "<value>"

# --- append new rows -------------------------------------------------------

# Row 32: SingleUseId39 -> "<value> l"
$ws.Range("B32").Value = "SingleUseId39"
$ws.Range("C32").Value = "Default"
$ws.Range("D32").Value = "Left"
$ws.Range("E32").Value = "LTR"
$ws.Range("F32").Value = "<value> l"

# Row 33: SingleUseId40 -> "Twoje zapotrzebowanie na wodę wynosi:"
$ws.Range("B33").Value = "SingleUseId40"
$ws.Range("C33").Value = "Default"
$ws.Range("D33").Value = "Left"
$ws.Range("E33").Value = "LTR"
$ws.Range("F33").Value = "Twoje zapotrzebowanie na wodę wynosi:"

# Row 34: SingleUseId41 -> "cm"
$ws.Range("B34").Value = "SingleUseId41"
$ws.Range("C34").Value = "Default"
$ws.Range("D34").Value = "Left"
$ws.Range("E34").Value = "LTR"
$ws.Range("F34").Value = "cm"

# Row 35: SingleUseId42 -> "l"
$ws.Range("B35").Value = "SingleUseId42"
$ws.Range("C35").Value = "Default"
$ws.Range("D35").Value = "Left"
$ws.Range("E35").Value = "LTR"
$ws.Range("F35").Value = "l"
